# Commit: "converter jin to kg unit"
# The 食材(ingredient) sheet recorded several "重量(公斤)" (weight) values
# using ad-hoc Taiwanese units (斤 "jin", 隻 "count", 大包 "big bag"). This
# converts those free-text quantities to plain decimal kilogram amounts
# (kept as text cells, same as the original, not numeric cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new "重量(公斤)" (column N) text value
$updates = @(
    @{ Addr = "N2";  Value = "15.0" },               # 25.8斤 -> 15.0
    @{ Addr = "N3";  Value = "7.199999999999999" },  # 100隻  -> 7.199999999999999
    @{ Addr = "N4";  Value = "0.6" },                 # 1大包  -> 0.6
    @{ Addr = "N5";  Value = "1.2" },                 # 2斤    -> 1.2
    @{ Addr = "N6";  Value = "1.2" },                 # 2斤    -> 1.2
    @{ Addr = "N7";  Value = "1.2" },                 # 2斤    -> 1.2
    @{ Addr = "N8";  Value = "4.8" },                 # 8斤    -> 4.8
    @{ Addr = "N9";  Value = "0.6" },                 # 1斤    -> 0.6
    @{ Addr = "N10"; Value = "2.4" }                  # 4斤    -> 2.4
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Addr)
    # Force a Text number format before assigning so Excel stores the
    # digit-looking value as a string (t="s") instead of auto-converting
    # it to a number, then clear the format again so the cell keeps its
    # original (default) style.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.ClearFormats()
}
